$d = $word.ActiveDocument
$wdFindContinue = 1
$wdReplaceAll = 2

function Replace-Text($search, $replace) {
    $d.Content.Find.Execute($search, $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, $replace, $wdReplaceAll) | Out-Null
}

# --- Paragraph: "Date of Meeting : 27/02/19" -- merge runs, drop proofErr ---
Replace-Text "Date of Meeting : " "Date of Meeting : "

# --- Paragraph: "Time of Meeting : 12:30" -- merge runs, drop proofErr ---
Replace-Text "Time of Meeting : 12:30" "Time of Meeting : 12:30"

# --- Paragraph: Attendees ---
Replace-Text "ees:- Joe, Andreea, Kacey" "ees:- Joe, Andreea, Kacey"

# --- Paragraph: Apologies ---
Replace-Text "Apologies from:- " "Apologies from:- "

# --- Paragraph: Item One ---
Replace-Text "Item One:-  Postmortem of previous week" "Item One:-  Postmortem of previous week"

# --- Paragraph: What went well ---
Replace-Text "What went well :  " "What went well :  "

# --- Paragraph: What went badly + new sentence ---
Replace-Text "What went badly : " "What went badly : "
Replace-Text "Some tasks were set after the start of the sprint." "Some tasks were set after the start of the sprint. Our Jira usage was not efficient, with us only creating one week’s worth of tasks at a time."

# --- Paragraph: Feedback Recieved ---
Replace-Text "Feedback Recieved : " "Feedback Recieved : "

# --- Paragraph: Individual work completed ---
Replace-Text "Individual work completed:-" "Individual work completed:-"

# --- Paragraph: Andreea research ---
Replace-Text "Andreea " "Andreea "

# --- Paragraph: Item 2 + new sentence ---
Replace-Text "Item 2:-  " "Item 2:-  "
Replace-Text "deliverable presentation for next week." "deliverable presentation for next week. We have created a long backlog of issues so we can more easily prepare sprints going forward."

# --- Paragraph: Tasks for the current week ---
Replace-Text "Tasks for the current week:-" "Tasks for the current week:-"

# --- Paragraph: Modify scripts ... bookmark removed from here, moved to Item 2 paragraph ---
Replace-Text "Modify scripts to allow mechanical change of direction" "Modify scripts to allow mechanical change of direction"

# Move the _GoBack bookmark from the "Modify scripts..." bullet paragraph to the
# end of the "Item 2..." paragraph (after the newly appended sentence).
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()

$rng = $d.Content
$rng.Find.Execute("going forward.", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "", 0) | Out-Null
$rng.Collapse(0)
$rng.InsertAfter([char]1)
$s = $rng.Start
$e = $rng.End
$bmRange = $d.Range($s, $e)
$d.Bookmarks.Add("_GoBack", $bmRange)
$d.Range($s, $e).Text = ""

# --- Paragraph: Meeting Ended ---
Replace-Text "Meeting Ended :- 13:00" "Meeting Ended :- 13:00"

# --- Paragraph: Minute Taker + Andreea ---
Replace-Text "Minute Taker:- Joe" "Minute Taker:- Joe"
$rng2 = $d.Content
$rng2.Find.Execute("Minute Taker:- Joe", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "", 0) | Out-Null
$rng2.Collapse(0)
$rng2.InsertAfter(" & Andreea")
